$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Action1")

# Update values: B2 "d" -> "d2", B3 "f" -> "f2"
$ws.Range("B2").Value = "d2"
$ws.Range("B3").Value = "f2"

# Move the active selection from B2 to B3
$ws.Range("B3").Select()
